# Add a "Timezone" column to the "Viajes" sheet.
#
# Before:  ... K=Aux_time_str  L=Parameter_avoid
# After:   ... K=Aux_time_str  L=Timezone          M=Parameter_avoid (new)
#
# The existing "Parameter_avoid" column (L) is pushed out to the new
# column M, and the freed-up column L becomes "Timezone". Row 2's data
# follows the same shift: the old L2 value ("tolls") moves to M2, the
# old K2 value ("73_148_ES") is cleared out (its column no longer holds
# data), and the new L2 is left blank too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Viajes")

# --- Row 1 (headers) ---------------------------------------------------

# M1 gets the header that used to live in L1 ("Parameter_avoid"), using
# the same bold/bordered header formatting as the rest of row 1.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$ws.Range("M1").Value = "Parameter_avoid"

# L1 becomes the new "Timezone" header. It keeps the plain (non-bold,
# non-bordered) style used by the other non-header-ish label cell (A2's
# style), matching the target formatting.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Range("L1").Value = "Timezone"

# --- Row 2 (sample data) ------------------------------------------------

# M2 gets the value that used to live in L2 ("tolls").
$ws.Range("A2").Copy() | Out-Null
$ws.Range("M2").PasteSpecial(-4122) | Out-Null
$ws.Range("M2").Value = "tolls"

# Old L2 ("tolls") and K2 ("73_148_ES") are now blank - the column they
# used to occupy doesn't carry data anymore. Give them the plain style
# (no more quote-prefixed text style on K2).
$ws.Range("K2").ClearContents()
$ws.Range("A2").Copy() | Out-Null
$ws.Range("K2").PasteSpecial(-4122) | Out-Null

$ws.Range("L2").ClearContents()
$ws.Range("A2").Copy() | Out-Null
$ws.Range("L2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
